# Update the wording on the "Methodology" SmartArt (ChevronBlockProcess)
# diagram on slide 2: make the chevron headline labels a bit more
# descriptive ("Mine and clean" -> "Mine and Clean", "Create" -> "Create
# Control", "Merge" -> "Merge Data", "Loop through" -> "Loop Through",
# "Conduct" -> "Conduct Test"). The longer description text under each
# chevron is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(5)
$sa = $shp.SmartArt

$sa.AllNodes.Item(1).TextFrame.TextRange.Text = "Mine and Clean"
$sa.AllNodes.Item(3).TextFrame.TextRange.Text = "Create Control"
$sa.AllNodes.Item(5).TextFrame.TextRange.Text = "Merge Data"
$sa.AllNodes.Item(7).TextFrame.TextRange.Text = "Loop Through"
$sa.AllNodes.Item(9).TextFrame.TextRange.Text = "Conduct Test"
